# "download articles with pandoc title blocks"
#
# Turns the old "Heading1 title + bold 'By <Author>' paragraph" (wrapped
# in a bookmark) into a pandoc-style title block: a paragraph styled
# "Title" holding the article title (one run per word/space/punctuation
# token) followed by a paragraph styled "Authors" holding just the
# author's name (again split word-by-word), with the enclosing bookmark
# removed and the leading "By " stripped.

$d = $word.ActiveDocument

$titlePara  = $d.Paragraphs.Item(1)   # "On Pilgrimage - December 1958" (Heading1, bookmarked)
$authorPara = $d.Paragraphs.Item(2)   # "By Dorothy Day" (bold)
$nextPara   = $d.Paragraphs.Item(3)   # "The Catholic Worker, December 1958, 1, 7."

# Remove both old paragraphs outright (this also collapses the
# bookmarkStart/bookmarkEnd pair that wraps the title paragraph down to
# a single zero-length point at the top of the document).
$victim = $d.Range($titlePara.Range.Start, $nextPara.Range.Start)
$victim.Delete()

# The now-orphaned bookmark markers sit at position 0; two zero-length
# deletes there sweep them out of the document entirely.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# Rebuild the title/author paragraphs as a pandoc-style title block,
# splitting every word, dash, and separating space into its own run.
$titleBlockXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">On</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">-</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">December</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">1958</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Range(0, 0).InsertXML($titleBlockXml)
